$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 18
$ws.Cells.Item(18, 7).Value2 = 2.7
$ws.Cells.Item(18, 9).Value2 = 2.88
$ws.Cells.Item(18, 12).Value2 = 4
$ws.Cells.Item(18, 24).Value2 = 11
$ws.Cells.Item(18, 27).Value2 = 29
$ws.Cells.Item(18, 35).Value2 = 12
$ws.Cells.Item(18, 36).Value2 = 13
$ws.Cells.Item(18, 37).Value2 = 34
$ws.Cells.Item(18, 38).Value2 = 34
$ws.Cells.Item(18, 40).Value2 = 4.5
$ws.Cells.Item(18, 41).Value2 = 17
$ws.Cells.Item(18, 42).Value2 = 34
$ws.Cells.Item(18, 51).Value2 = 4.75
$ws.Cells.Item(18, 53).Value2 = 41

# Row 19
$ws.Cells.Item(19, 7).Value2 = 1.9
$ws.Cells.Item(19, 10).Value2 = 2.75
$ws.Cells.Item(19, 12).Value2 = 5.5
$ws.Cells.Item(19, 23).Value2 = 5
$ws.Cells.Item(19, 26).Value2 = 15
$ws.Cells.Item(19, 30).Value2 = 6.5
$ws.Cells.Item(19, 31).Value2 = 23
$ws.Cells.Item(19, 34).Value2 = 8.5
$ws.Cells.Item(19, 35).Value2 = 21
$ws.Cells.Item(19, 36).Value2 = 17
$ws.Cells.Item(19, 38).Value2 = 41
$ws.Cells.Item(19, 39).Value2 = 51
$ws.Cells.Item(19, 40).Value2 = 3.6
$ws.Cells.Item(19, 41).Value2 = 11
$ws.Cells.Item(19, 52).Value2 = 29
$ws.Cells.Item(19, 54).Value2 = 126

# Row 23
$ws.Cells.Item(23, 7).Value2 = 3.5
$ws.Cells.Item(23, 9).Value2 = 2.35
$ws.Cells.Item(23, 10).Value2 = 4
$ws.Cells.Item(23, 15).Value2 = 1.5
$ws.Cells.Item(23, 16).Value2 = 2.5
$ws.Cells.Item(23, 17).Value2 = 2.6
$ws.Cells.Item(23, 18).Value2 = 1.48
$ws.Cells.Item(23, 19).Value2 = 1.57
$ws.Cells.Item(23, 20).Value2 = 2.25
$ws.Cells.Item(23, 22).Value2 = 1.62
$ws.Cells.Item(23, 32).Value2 = 67
$ws.Cells.Item(23, 35).Value2 = 10
$ws.Cells.Item(23, 43).Value2 = 67
$ws.Cells.Item(23, 46).Value2 = 2.25
$ws.Cells.Item(23, 51).Value2 = 4.33

# Row 45
$ws.Cells.Item(45, 7).Value2 = 1.5
$ws.Cells.Item(45, 8).Value2 = 4.33
$ws.Cells.Item(45, 9).Value2 = 6.5
$ws.Cells.Item(45, 10).Value2 = 2.05
$ws.Cells.Item(45, 15).Value2 = 1.22
$ws.Cells.Item(45, 16).Value2 = 4
$ws.Cells.Item(45, 17).Value2 = 1.73
$ws.Cells.Item(45, 18).Value2 = 2.08
$ws.Cells.Item(45, 19).Value2 = 1.33
$ws.Cells.Item(45, 20).Value2 = 3.25
$ws.Cells.Item(45, 25).Value2 = 8.5
$ws.Cells.Item(45, 26).Value2 = 11
$ws.Cells.Item(45, 30).Value2 = 8
$ws.Cells.Item(45, 34).Value2 = 17
$ws.Cells.Item(45, 36).Value2 = 19
$ws.Cells.Item(45, 41).Value2 = 7.5
$ws.Cells.Item(45, 46).Value2 = 3.25
$ws.Cells.Item(45, 52).Value2 = 29

# Row 46
$ws.Cells.Item(46, 7).Value2 = 1.47
$ws.Cells.Item(46, 8).Value2 = 4.15
$ws.Cells.Item(46, 9).Value2 = 6
$ws.Cells.Item(46, 10).Value2 = 1.98
$ws.Cells.Item(46, 11).Value2 = 2.27
$ws.Cells.Item(46, 12).Value2 = 5.7
$ws.Cells.Item(46, 13).Value2 = 1.03
$ws.Cells.Item(46, 14).Value2 = 11.5
$ws.Cells.Item(46, 17).Value2 = 1.7
$ws.Cells.Item(46, 18).Value2 = 1.91
$ws.Cells.Item(46, 20).Value2 = 3.13
$ws.Cells.Item(46, 21).Value2 = 1.87
$ws.Cells.Item(46, 22).Value2 = 1.75
$ws.Cells.Item(46, 24).Value2 = 6.8
$ws.Cells.Item(46, 26).Value2 = 9.75
$ws.Cells.Item(46, 29).Value2 = 11.5
$ws.Cells.Item(46, 30).Value2 = 8.25
$ws.Cells.Item(46, 31).Value2 = 18.5
$ws.Cells.Item(46, 32).Value2 = 90
$ws.Cells.Item(46, 34).Value2 = 16
$ws.Cells.Item(46, 35).Value2 = 37
$ws.Cells.Item(46, 36).Value2 = 19
$ws.Cells.Item(46, 37).Value2 = 120
$ws.Cells.Item(46, 38).Value2 = 65
$ws.Cells.Item(46, 39).Value2 = 65
$ws.Cells.Item(46, 40).Value2 = 3.25
$ws.Cells.Item(46, 41).Value2 = 6.8
$ws.Cells.Item(46, 43).Value2 = 20
$ws.Cells.Item(46, 44).Value2 = 50
$ws.Cells.Item(46, 46).Value2 = 2.85
$ws.Cells.Item(46, 47).Value2 = 8.25
$ws.Cells.Item(46, 48).Value2 = 80
$ws.Cells.Item(46, 51).Value2 = 7.2
$ws.Cells.Item(46, 52).Value2 = 35
$ws.Cells.Item(46, 53).Value2 = 37
$ws.Cells.Item(46, 54).Value2 = 250

# Row 47
$ws.Cells.Item(47, 7).Value2 = 1.35
$ws.Cells.Item(47, 8).Value2 = 4.55
$ws.Cells.Item(47, 9).Value2 = 7.8
$ws.Cells.Item(47, 10).Value2 = 1.82
$ws.Cells.Item(47, 11).Value2 = 2.4
$ws.Cells.Item(47, 12).Value2 = 6.6
$ws.Cells.Item(47, 14).Value2 = 13
$ws.Cells.Item(47, 19).Value2 = 1.28
$ws.Cells.Item(47, 20).Value2 = 3.54
$ws.Cells.Item(47, 22).Value2 = 1.8
$ws.Cells.Item(47, 23).Value2 = 7.7
$ws.Cells.Item(47, 24).Value2 = 6.9
$ws.Cells.Item(47, 29).Value2 = 14
$ws.Cells.Item(47, 34).Value2 = 23
$ws.Cells.Item(47, 35).Value2 = 60
$ws.Cells.Item(47, 36).Value2 = 24
$ws.Cells.Item(47, 37).Value2 = 200
$ws.Cells.Item(47, 38).Value2 = 90
$ws.Cells.Item(47, 39).Value2 = 65
$ws.Cells.Item(47, 40).Value2 = 3.2
$ws.Cells.Item(47, 41).Value2 = 6
$ws.Cells.Item(47, 42).Value2 = 15
$ws.Cells.Item(47, 43).Value2 = 16
$ws.Cells.Item(47, 45).Value2 = 175
$ws.Cells.Item(47, 47).Value2 = 7.9
$ws.Cells.Item(47, 48).Value2 = 65
$ws.Cells.Item(47, 51).Value2 = 8.75
$ws.Cells.Item(47, 52).Value2 = 45
$ws.Cells.Item(47, 53).Value2 = 40
$ws.Cells.Item(47, 54).Value2 = 300
$ws.Cells.Item(47, 55).Value2 = 300

# Row 86
$ws.Cells.Item(86, 7).Value2 = 1.55
$ws.Cells.Item(86, 8).Value2 = 3.9
$ws.Cells.Item(86, 9).Value2 = 5.1
$ws.Cells.Item(86, 10).Value2 = 2.12
$ws.Cells.Item(86, 11).Value2 = 2.22
$ws.Cells.Item(86, 12).Value2 = 5.4
$ws.Cells.Item(86, 14).Value2 = 7.7
$ws.Cells.Item(86, 15).Value2 = 1.27
$ws.Cells.Item(86, 16).Value2 = 3.4
$ws.Cells.Item(86, 17).Value2 = 1.82
$ws.Cells.Item(86, 18).Value2 = 1.93
$ws.Cells.Item(86, 21).Value2 = 1.88
$ws.Cells.Item(86, 22).Value2 = 1.82
$ws.Cells.Item(86, 23).Value2 = 6.8
$ws.Cells.Item(86, 24).Value2 = 7.2
$ws.Cells.Item(86, 26).Value2 = 11
$ws.Cells.Item(86, 27).Value2 = 12.5
$ws.Cells.Item(86, 29).Value2 = 7.7
$ws.Cells.Item(86, 30).Value2 = 7.7
$ws.Cells.Item(86, 31).Value2 = 17.5
$ws.Cells.Item(86, 32).Value2 = 90
$ws.Cells.Item(86, 34).Value2 = 14
$ws.Cells.Item(86, 35).Value2 = 30
$ws.Cells.Item(86, 36).Value2 = 17
$ws.Cells.Item(86, 37).Value2 = 100
$ws.Cells.Item(86, 38).Value2 = 55
$ws.Cells.Item(86, 39).Value2 = 55
$ws.Cells.Item(86, 40).Value2 = 3.35
$ws.Cells.Item(86, 41).Value2 = 7.5
$ws.Cells.Item(86, 42).Value2 = 18
$ws.Cells.Item(86, 43).Value2 = 24
$ws.Cells.Item(86, 47).Value2 = 8.25
$ws.Cells.Item(86, 48).Value2 = 80
$ws.Cells.Item(86, 51).Value2 = 6.8
$ws.Cells.Item(86, 52).Value2 = 30
$ws.Cells.Item(86, 53).Value2 = 37
$ws.Cells.Item(86, 54).Value2 = 200
$ws.Cells.Item(86, 55).Value2 = 250
$ws.Cells.Item(86, 56).Value2 = 500
